$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C) column values from 45192 to 45202 for existing rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 10 gains an explicit 15pt custom row height
$ws.Rows.Item(10).RowHeight = 15

# Add new data row 11
$ws.Range("A11").Value = "A 45983-2023"
$ws.Range("B11").Value = 45196
$ws.Range("C11").Value = 45202
$ws.Range("B11:C11").NumberFormat = "YYYY-MM-DD"
$ws.Range("D11").Value = "SKÅNE LÄN"
$ws.Range("E11").Value = "SKURUP"
$ws.Range("G11").Value = 0.6
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0

# R11 mirrors the empty, wrap-text styled placeholder cells used in R2:R10
$ws.Range("R10").Copy($ws.Range("R11"))
$ws.Range("R11").Formula = ""
